$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "November" data row (row 3).
$ws.Range("A3").Value = 2024
$ws.Range("B3").Value = 11
$ws.Range("C3").Value = "November"
$ws.Range("D3").Value = 11
$ws.Range("E3").Formula = "=D3+6"
$ws.Range("F3").Value = 7

# Rename header labels for start/end columns (D1, E1) to the new "Day" wording.
$ws.Range("D1").Value = "startDay"
$ws.Range("E1").Value = "endDay"

# Match the saved selection/active cell from the target workbook.
$ws.Range("D2").Select()
